$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6 and 7 swap their ENVO/MOP IRI + desc content (CHEBI_52214 <-> CHEBI_39141)
$ws.Range("B6").Value = "http://purl.obolibrary.org/obo/CHEBI_39141"
$ws.Range("C6").Value = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39141'}"
$ws.Range("D6").Value = "http://purl.obolibrary.org/obo/CHEBI_39141"
$ws.Range("E6").Value = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39141'}"

$ws.Range("B7").Value = "http://purl.obolibrary.org/obo/CHEBI_52214"
$ws.Range("C7").Value = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_52214'}"
$ws.Range("D7").Value = "http://purl.obolibrary.org/obo/CHEBI_52214"
$ws.Range("E7").Value = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_52214'}"

# Add new header "MOP_DEF" in F1, copying formatting from E1 (bold, bordered header style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "MOP_DEF"

# Populate MOP_DEF column values for each data row (2-26)
$ws.Range("F2").Value = "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(`"Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'`", 'en')]"
$ws.Range("F3").Value = "['B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]']"
$ws.Range("F4").Value = "['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']"
$ws.Range("F5").Value = "['A subatomic particle is a material that is below the scale of an atom. [Allotrope]']"

for ($r = 6; $r -le 26; $r++) {
    $ws.Range("F$r").Value = "[]"
}
